# Update the "Codebase copy & modification" caption textbox (slide 1) to
# read "Copied &" / "modified" on two separate paragraphs. PowerPoint's
# spAutoFit will shrink the shape's height to match the new two-line text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(82)   # "TextBox 96"

$sh.TextFrame.TextRange.Text = "Copied &" + [char]13 + "modified"
